# Weekly data refresh: insert two new rows of "Membrillo" price records
# (fecha 2023-04-28 / serial 45044) above the existing history, pushing the
# prior rows 7-11 down to rows 9-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 7; everything below (old
# rows 7-11) shifts down to rows 9-13, carrying its values/formatting along.
$ws.Rows("7:8").Insert()

# New row 7: Especial quality, 18kg empedrada box, Región de O'Higgins
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 45044
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100104
$ws.Range("H7").Value = "Frutos de pepita"
$ws.Range("I7").Value = 100104003
$ws.Range("J7").Value = "Membrillo"
$ws.Range("K7").Value = "Champion"
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range("Q7").Value = "$/caja 18 kilos empedrada"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 722
$ws.Range("T7").Value = 18

# New row 8: Primera quality, 18kg empedrada box, Región de O'Higgins
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 45044
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100104
$ws.Range("H8").Value = "Frutos de pepita"
$ws.Range("I8").Value = 100104003
$ws.Range("J8").Value = "Membrillo"
$ws.Range("K8").Value = "Champion"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("Q8").Value = "$/caja 18 kilos empedrada"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 667
$ws.Range("T8").Value = 18
